# Import dosyalarında excelden alınacak şekilde ZeyilSirketTanimli tablosu eklendi.
$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it at the end.
$originalActive = $wb.Worksheets.Item(1)

# --- Update the selection remembered on the "SirketTanimliUrun" sheet ---
$ws3 = $wb.Worksheets.Item("SirketTanimliUrun")
$ws3.Activate() | Out-Null
$ws3.Range("B30").Select() | Out-Null

# --- Add the new "Zeyiller" sheet at the end of the workbook ---
$ws4 = $wb.Worksheets.Add()
$ws4.Name = "Zeyiller"
$ws4.Range("A1").Value = "ZeyilKod"
$ws4.Range("A2").Select() | Out-Null

# Move it after the last existing sheet so it lands at the end.
$ws4.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Restore the originally active sheet.
$originalActive.Activate() | Out-Null
